$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 (Test_details): row 16 no longer needs its explicit 30pt height ---
$ws1.Rows.Item(16).EntireRow.AutoFit()

# --- sheet2 (Test_scenarios_Cart_page): add TC005 scenario row ---
# The text that used to describe TC004 ("search for an item...") is being
# moved down to a new row and relabelled TC005; TC004 gets a new description
# about the 'Send Us Your Comments' flow.
$existingDescription = $ws2.Range("E6").Value2

$ws2.Range("D7").Value = "TC005"
$ws2.Range("E7").Value = $existingDescription
$ws2.Range("E6").Value = "For an unlogged user, open Cart page on its own, click on the 'Send Us Your Comments' link, verify if it brings user to the correct page, and fill in a feedback form"

# Copy formatting (borders/alignment/wrap) from row 6 down into row 7 so the
# new row matches the rest of the table
$ws2.Range("B6:E6").Copy()
$ws2.Range("B7:E7").PasteSpecial(-4122)

# Row 6's text got shorter (needs less height); row 7 carries the longer text
# that used to live in row 6
$ws2.Rows.Item(6).RowHeight = 60
$ws2.Rows.Item(7).RowHeight = 75

# Leave the selection on the newly-added description cell ...
$ws2.Range("E7").Select()

# ... but keep Test_details as the active/visible tab
$ws1.Activate()

$wb.Save()
